$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 84400
$ws.Range("D2").Value = 77500
$ws.Range("C3").Value = 46400
$ws.Range("D3").Value = 68100
$ws.Range("C4").Value = 125000
$ws.Range("D4").Value = -30000
$ws.Range("C5").Value = 110000
$ws.Range("D5").Value = -30000
$ws.Range("C6").Value = 6600
$ws.Range("D6").Value = 62100
$ws.Range("C7").Value = 93200
$ws.Range("D7").Value = -30000
$ws.Range("D8").Value = 92000
$ws.Range("C9").Value = 90800
$ws.Range("D9").Value = -30000
$ws.Range("B10").Value = "Chia lester"
$ws.Range("C10").Value = 64500
$ws.Range("D10").Value = -25000
$ws.Range("C11").Value = 54600
$ws.Range("D11").Value = -30000
$ws.Range("C12").Value = 14800
$ws.Range("D12").Value = 6200
$ws.Range("C13").Value = 44600
$ws.Range("D13").Value = -30000
$ws.Range("D14").Value = 44000
$ws.Range("C15").Value = 7000
$ws.Range("D15").Value = 4000
$ws.Range("C16").Value = 400
$ws.Range("D16").Value = 4000
$ws.Range("C17").Value = -12900
$ws.Range("D17").Value = 14900
$ws.Range("D18").Value = 28800
$ws.Range("D19").Value = 28000
$ws.Range("D20").Value = 23900
$ws.Range("D21").Value = 17700
$ws.Range("C22").Value = 2400
$ws.Range("D22").Value = -30000
$ws.Range("C23").Value = -5400
$ws.Range("D23").Value = -30000
$ws.Range("C24").Value = -10200
$ws.Range("D24").Value = -30000
$ws.Range("C25").Value = -16800
$ws.Range("D25").Value = -30000
$ws.Range("C42").Value = -30000
$ws.Range("D42").Value = -30000

$ws.Range("A2:A45").Select()
